# Updated cryptos list on Mon Aug  7 22:40:44 UTC 2023 with GitHub Actions
# Refresh Price (D) / Volume(1h) (E) columns, and swap two reordered rows B/C text
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '29.140.85'
$ws.Range("E2").Value = '  +0.09%  '

$ws.Range("D3").Value = '1.823.21'

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '0.9983'
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = '  -0.15%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '241.29'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -0.83%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '0.6200'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  -1.15%  '

$ws.Range("E7").Value = '  -0.14%  '

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.07340'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  -2.16%  '

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.2899'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  -1.11%  '

$ws.Range("E10").Value = '  -1.37%  '

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.07673'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  -0.37%  '

$ws.Range("D12").Value = '1.827.11'
$ws.Range("E12").Value = '  +0.85%  '

$ws.Range("E13").Value = '  -1.53%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '0.6633'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  -0.91%  '

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '82.22'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  -0.80%  '

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '0.000008954'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  -4.50%  '

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '5.827'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  -2.84%  '

$ws.Range("D18").Value = '29.121.58'
$ws.Range("E18").Value = '  +0.06%  '

$ws.Range("D19").Value = '2.068.64'
$ws.Range("E19").Value = '  +0.26%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '239.11'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  +7.05%  '

$ws.Range("E21").Value = '  -1.73%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '0.9993'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  -0.23%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '7.214'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  +0.99%  '

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '0.9997'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  -0.13%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '157.93'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  -1.43%  '

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '0.1423'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  +1.55%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '8.467'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  -0.59%  '

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '17.64'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  -1.51%  '

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '1.484'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  -0.49%  '

$ws.Range("E30").Value = '  -4.47%  '

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '4.083'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  -1.10%  '

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '4.091'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  -1.83%  '

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '1.204'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  -0.55%  '

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '1.835'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  +0.22%  '

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '0.7328'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  -1.32%  '

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '1.135'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  -0.48%  '

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '2.625'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  -1.71%  '

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '2.838'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  +2.75%  '

$ws.Range("D39").Value = '1.213.14'
$ws.Range("E39").Value = '  -1.78%  '

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.01762'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  -0.97%  '

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '6.297'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  -3.15%  '

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.9128'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  +2.10%  '

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.9997'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  -0.13%  '

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '101.59'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  -0.61%  '

$ws.Range("D45").Value = '1.975.59'
$ws.Range("E45").Value = '  +0.07%  '

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '64.60'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  -2.01%  '

$ws.Range("E47").Value = '  -0.14%  '

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '0.00000000118'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  -7.20%  '

$ws.Range("B49").Value = 'EnergySwap'
$ws.Range("C49").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '9.104'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  +1.01%  '

$ws.Range("B50").Value = 'TheSandbox'
$ws.Range("C50").Value = 'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand'
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.4008'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  -1.63%  '

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.05752'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  -1.32%  '
